# Quick registration module: Names Field (Done)
# Appends the newly-registered user (documento 66000093) as a new row
# at the bottom of the "usuarios_registrados" table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row right after the current table (row 52 here).
$newRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

$ws.Cells.Item($newRow, 1).Value = 66000093
$ws.Cells.Item($newRow, 2).Value = "66000093test@gmail.com"
$ws.Cells.Item($newRow, 3).Value = 66000093
$ws.Cells.Item($newRow, 4).Value = "Aaaaaaaaa1"
